$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(97, 4).Value = 44518
$ws.Cells.Item(97, 11).Value = 500
$ws.Cells.Item(97, 12).Value = 500
$ws.Cells.Item(97, 13).Value = 500
$ws.Cells.Item(97, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(97, 16).Value = 500
$ws.Cells.Item(98, 4).Value = 44518
$ws.Cells.Item(98, 8).Value = "Camote"
$ws.Cells.Item(98, 9).Value = "1a nueva(o)"
$ws.Cells.Item(98, 10).Value = 900
$ws.Cells.Item(98, 11).Value = 700
$ws.Cells.Item(98, 12).Value = 700
$ws.Cells.Item(98, 13).Value = 700
$ws.Cells.Item(98, 15).Value = "Perú"
$ws.Cells.Item(98, 16).Value = 700
$ws.Cells.Item(99, 4).Value = 44518
$ws.Cells.Item(99, 8).Value = "Paine"
$ws.Cells.Item(99, 9).Value = "1a (guarda)"
$ws.Cells.Item(99, 10).Value = 1500
$ws.Cells.Item(99, 11).Value = 80
$ws.Cells.Item(99, 12).Value = 80
$ws.Cells.Item(99, 13).Value = 80
$ws.Cells.Item(99, 16).Value = 80
$ws.Cells.Item(100, 4).Value = 44463
$ws.Cells.Item(100, 8).Value = "Camote"
$ws.Cells.Item(100, 10).Value = 800
$ws.Cells.Item(100, 11).Value = 400
$ws.Cells.Item(100, 12).Value = 400
$ws.Cells.Item(100, 13).Value = 400
$ws.Cells.Item(100, 16).Value = 400
$ws.Cells.Item(101, 4).Value = 44463
$ws.Cells.Item(101, 8).Value = "Paine"
$ws.Cells.Item(101, 9).Value = "1a (guarda)"
$ws.Cells.Item(101, 10).Value = 1500
$ws.Cells.Item(101, 11).Value = 130
$ws.Cells.Item(101, 12).Value = 130
$ws.Cells.Item(101, 13).Value = 130
$ws.Cells.Item(101, 16).Value = 130
$ws.Cells.Item(102, 4).Value = 44245
$ws.Cells.Item(102, 9).Value = "1a nueva(o)"
$ws.Cells.Item(102, 11).Value = 280
$ws.Cells.Item(102, 12).Value = 280
$ws.Cells.Item(102, 13).Value = 280
$ws.Cells.Item(102, 16).Value = 280
$ws.Cells.Item(103, 4).Value = 44481
$ws.Cells.Item(103, 11).Value = 100
$ws.Cells.Item(103, 12).Value = 120
$ws.Cells.Item(103, 13).Value = 110
$ws.Cells.Item(103, 16).Value = 110
$ws.Cells.Item(104, 4).Value = 44229
$ws.Cells.Item(104, 10).Value = 900
$ws.Cells.Item(105, 4).Value = 44417
$ws.Cells.Item(105, 10).Value = 800
$ws.Cells.Item(105, 11).Value = 450
$ws.Cells.Item(105, 12).Value = 450
$ws.Cells.Item(105, 13).Value = 450
$ws.Cells.Item(105, 16).Value = 450
$ws.Cells.Item(106, 4).Value = 44445
$ws.Cells.Item(106, 10).Value = 2000
$ws.Cells.Item(107, 4).Value = 44249
$ws.Cells.Item(107, 10).Value = 800
$ws.Cells.Item(107, 11).Value = 300
$ws.Cells.Item(107, 12).Value = 300
$ws.Cells.Item(107, 13).Value = 300
$ws.Cells.Item(107, 16).Value = 300
$ws.Cells.Item(108, 4).Value = 44342
$ws.Cells.Item(108, 9).Value = "1a (guarda)"
$ws.Cells.Item(108, 10).Value = 900
$ws.Cells.Item(108, 11).Value = 250
$ws.Cells.Item(108, 12).Value = 250
$ws.Cells.Item(108, 13).Value = 250
$ws.Cells.Item(108, 16).Value = 250
$ws.Cells.Item(109, 4).Value = 44342
$ws.Cells.Item(109, 8).Value = "Paine"
$ws.Cells.Item(109, 10).Value = 1200
$ws.Cells.Item(109, 11).Value = 150
$ws.Cells.Item(109, 12).Value = 150
$ws.Cells.Item(109, 13).Value = 150
$ws.Cells.Item(109, 16).Value = 150
$ws.Cells.Item(110, 4).Value = 44216
$ws.Cells.Item(110, 8).Value = "Camote"
$ws.Cells.Item(110, 9).Value = "1a nueva(o)"
$ws.Cells.Item(110, 10).Value = 1000
$ws.Cells.Item(110, 11).Value = 320
$ws.Cells.Item(110, 12).Value = 350
$ws.Cells.Item(110, 13).Value = 335
$ws.Cells.Item(110, 16).Value = 335
$ws.Cells.Item(111, 4).Value = 44270
$ws.Cells.Item(111, 9).Value = "1a (cosecha)"
$ws.Cells.Item(111, 10).Value = 1000
$ws.Cells.Item(111, 11).Value = 280
$ws.Cells.Item(111, 12).Value = 300
$ws.Cells.Item(111, 13).Value = 290
$ws.Cells.Item(111, 16).Value = 290
$ws.Cells.Item(112, 4).Value = 44363
$ws.Cells.Item(112, 9).Value = "1a (guarda)"
$ws.Cells.Item(112, 10).Value = 800
$ws.Cells.Item(112, 11).Value = 280
$ws.Cells.Item(112, 12).Value = 280
$ws.Cells.Item(112, 13).Value = 280
$ws.Cells.Item(112, 16).Value = 280
$ws.Cells.Item(113, 4).Value = 44363
$ws.Cells.Item(113, 8).Value = "Paine"
$ws.Cells.Item(113, 10).Value = 1500
$ws.Cells.Item(113, 11).Value = 150
$ws.Cells.Item(113, 12).Value = 150
$ws.Cells.Item(113, 13).Value = 150
$ws.Cells.Item(113, 16).Value = 150
$ws.Cells.Item(114, 4).Value = 44299
$ws.Cells.Item(114, 8).Value = "Camote"
$ws.Cells.Item(114, 10).Value = 900
$ws.Cells.Item(114, 11).Value = 200
$ws.Cells.Item(114, 12).Value = 200
$ws.Cells.Item(114, 13).Value = 200
$ws.Cells.Item(114, 16).Value = 200
$ws.Cells.Item(115, 4).Value = 44257
$ws.Cells.Item(115, 10).Value = 1000
$ws.Cells.Item(115, 11).Value = 300
$ws.Cells.Item(115, 12).Value = 300
$ws.Cells.Item(115, 13).Value = 300
$ws.Cells.Item(115, 16).Value = 300
$ws.Cells.Item(116, 4).Value = 44372
$ws.Cells.Item(116, 9).Value = "1a (guarda)"
$ws.Cells.Item(116, 10).Value = 1000
$ws.Cells.Item(116, 11).Value = 200
$ws.Cells.Item(116, 12).Value = 200
$ws.Cells.Item(116, 13).Value = 200
$ws.Cells.Item(116, 16).Value = 200
$ws.Cells.Item(117, 4).Value = 44372
$ws.Cells.Item(117, 8).Value = "Paine"
$ws.Cells.Item(117, 10).Value = 1500
$ws.Cells.Item(117, 11).Value = 140
$ws.Cells.Item(117, 12).Value = 140
$ws.Cells.Item(117, 13).Value = 140
$ws.Cells.Item(117, 16).Value = 140
$ws.Cells.Item(118, 4).Value = 44169
$ws.Cells.Item(118, 11).Value = 1000
$ws.Cells.Item(118, 12).Value = 1000
$ws.Cells.Item(118, 13).Value = 1000
$ws.Cells.Item(118, 16).Value = 1000
$ws.Cells.Item(119, 4).Value = 44195
$ws.Cells.Item(119, 9).Value = "1a nueva(o)"
$ws.Cells.Item(119, 11).Value = 400
$ws.Cells.Item(119, 12).Value = 400
$ws.Cells.Item(119, 13).Value = 400
$ws.Cells.Item(119, 16).Value = 400
$ws.Cells.Item(120, 4).Value = 44376
$ws.Cells.Item(120, 9).Value = "1a (guarda)"
$ws.Cells.Item(120, 10).Value = 900
$ws.Cells.Item(120, 11).Value = 200
$ws.Cells.Item(120, 12).Value = 200
$ws.Cells.Item(120, 13).Value = 200
$ws.Cells.Item(120, 16).Value = 200
$ws.Cells.Item(121, 4).Value = 44242
$ws.Cells.Item(121, 10).Value = 800
$ws.Cells.Item(121, 11).Value = 300
$ws.Cells.Item(121, 12).Value = 300
$ws.Cells.Item(121, 13).Value = 300
$ws.Cells.Item(121, 16).Value = 300
$ws.Cells.Item(122, 4).Value = 44431
$ws.Cells.Item(122, 11).Value = 600
$ws.Cells.Item(122, 12).Value = 600
$ws.Cells.Item(122, 13).Value = 600
$ws.Cells.Item(122, 16).Value = 600
$ws.Cells.Item(123, 4).Value = 44239
$ws.Cells.Item(123, 10).Value = 800
$ws.Cells.Item(123, 11).Value = 300
$ws.Cells.Item(123, 12).Value = 300
$ws.Cells.Item(123, 13).Value = 300
$ws.Cells.Item(123, 16).Value = 300
$ws.Cells.Item(124, 4).Value = 44222
$ws.Cells.Item(124, 9).Value = "1a nueva(o)"
$ws.Cells.Item(124, 10).Value = 1000
$ws.Cells.Item(124, 11).Value = 320
$ws.Cells.Item(124, 12).Value = 350
$ws.Cells.Item(124, 13).Value = 335
$ws.Cells.Item(124, 16).Value = 335
$ws.Cells.Item(125, 4).Value = 44426
$ws.Cells.Item(125, 10).Value = 800
$ws.Cells.Item(125, 11).Value = 430
$ws.Cells.Item(125, 12).Value = 430
$ws.Cells.Item(125, 13).Value = 430
$ws.Cells.Item(125, 16).Value = 430
$ws.Cells.Item(126, 4).Value = 44176
$ws.Cells.Item(126, 9).Value = "1a nueva(o)"
$ws.Cells.Item(126, 11).Value = 850
$ws.Cells.Item(126, 12).Value = 850
$ws.Cells.Item(126, 13).Value = 850
$ws.Cells.Item(126, 16).Value = 850
$ws.Cells.Item(127, 4).Value = 44284
$ws.Cells.Item(127, 8).Value = "Camote"
$ws.Cells.Item(127, 9).Value = "1a (cosecha)"
$ws.Cells.Item(127, 10).Value = 800
$ws.Cells.Item(127, 11).Value = 200
$ws.Cells.Item(127, 12).Value = 200
$ws.Cells.Item(127, 13).Value = 200
$ws.Cells.Item(127, 16).Value = 200
$ws.Cells.Item(128, 4).Value = 44441
$ws.Cells.Item(128, 11).Value = 800
$ws.Cells.Item(128, 12).Value = 800
$ws.Cells.Item(128, 13).Value = 800
$ws.Cells.Item(128, 16).Value = 800
$ws.Cells.Item(129, 4).Value = 44279
$ws.Cells.Item(129, 8).Value = "Camote"
$ws.Cells.Item(129, 9).Value = "1a (cosecha)"
$ws.Cells.Item(129, 10).Value = 900
$ws.Cells.Item(129, 11).Value = 200
$ws.Cells.Item(129, 12).Value = 200
$ws.Cells.Item(129, 13).Value = 200
$ws.Cells.Item(129, 16).Value = 200
$ws.Cells.Item(130, 4).Value = 44504
$ws.Cells.Item(130, 8).Value = "Paine"
$ws.Cells.Item(130, 10).Value = 2500
$ws.Cells.Item(130, 11).Value = 80
$ws.Cells.Item(130, 12).Value = 80
$ws.Cells.Item(130, 13).Value = 80
$ws.Cells.Item(130, 16).Value = 80
$ws.Cells.Item(131, 4).Value = 44350
$ws.Cells.Item(131, 10).Value = 900
$ws.Cells.Item(131, 11).Value = 280
$ws.Cells.Item(131, 12).Value = 280
$ws.Cells.Item(131, 13).Value = 280
$ws.Cells.Item(131, 16).Value = 280
$ws.Cells.Item(132, 4).Value = 44350
$ws.Cells.Item(132, 8).Value = "Paine"
$ws.Cells.Item(132, 9).Value = "1a (guarda)"
$ws.Cells.Item(132, 10).Value = 1200
$ws.Cells.Item(132, 11).Value = 150
$ws.Cells.Item(132, 12).Value = 150
$ws.Cells.Item(132, 13).Value = 150
$ws.Cells.Item(132, 16).Value = 150
$ws.Cells.Item(133, 4).Value = 44384
$ws.Cells.Item(133, 8).Value = "Camote"
$ws.Cells.Item(133, 9).Value = "1a (guarda)"
$ws.Cells.Item(133, 10).Value = 600
$ws.Cells.Item(133, 11).Value = 300
$ws.Cells.Item(133, 12).Value = 300
$ws.Cells.Item(133, 13).Value = 300
$ws.Cells.Item(133, 16).Value = 300
$ws.Cells.Item(134, 4).Value = 44329
$ws.Cells.Item(134, 10).Value = 800
$ws.Cells.Item(134, 11).Value = 300
$ws.Cells.Item(134, 12).Value = 300
$ws.Cells.Item(134, 13).Value = 300
$ws.Cells.Item(134, 16).Value = 300
$ws.Cells.Item(135, 4).Value = 44272
$ws.Cells.Item(135, 9).Value = "1a (cosecha)"
$ws.Cells.Item(135, 10).Value = 900
$ws.Cells.Item(135, 11).Value = 250
$ws.Cells.Item(135, 12).Value = 250
$ws.Cells.Item(135, 13).Value = 250
$ws.Cells.Item(135, 16).Value = 250
$ws.Cells.Item(136, 4).Value = 44272
$ws.Cells.Item(136, 8).Value = "Paine"
$ws.Cells.Item(136, 9).Value = "1a (cosecha)"
$ws.Cells.Item(136, 10).Value = 1200
$ws.Cells.Item(136, 11).Value = 150
$ws.Cells.Item(136, 12).Value = 150
$ws.Cells.Item(136, 13).Value = 150
$ws.Cells.Item(136, 16).Value = 150
$ws.Cells.Item(137, 4).Value = 44305
$ws.Cells.Item(137, 8).Value = "Camote"
$ws.Cells.Item(137, 10).Value = 900
$ws.Cells.Item(137, 11).Value = 200
$ws.Cells.Item(137, 12).Value = 200
$ws.Cells.Item(137, 13).Value = 200
$ws.Cells.Item(137, 16).Value = 200
$ws.Cells.Item(138, 4).Value = 44225
$ws.Cells.Item(138, 9).Value = "1a nueva(o)"
$ws.Cells.Item(138, 11).Value = 300
$ws.Cells.Item(138, 12).Value = 300
$ws.Cells.Item(138, 13).Value = 300
$ws.Cells.Item(138, 16).Value = 300
$ws.Cells.Item(139, 4).Value = 44348
$ws.Cells.Item(139, 9).Value = "1a (guarda)"
$ws.Cells.Item(139, 10).Value = 900
$ws.Cells.Item(139, 11).Value = 280
$ws.Cells.Item(139, 12).Value = 280
$ws.Cells.Item(139, 13).Value = 280
$ws.Cells.Item(139, 16).Value = 280
$ws.Cells.Item(140, 4).Value = 44348
$ws.Cells.Item(140, 8).Value = "Paine"
$ws.Cells.Item(140, 10).Value = 1500
$ws.Cells.Item(140, 11).Value = 150
$ws.Cells.Item(140, 12).Value = 150
$ws.Cells.Item(140, 13).Value = 150
$ws.Cells.Item(140, 16).Value = 150
$ws.Cells.Item(141, 4).Value = 44322
$ws.Cells.Item(141, 8).Value = "Camote"
$ws.Cells.Item(141, 10).Value = 800
$ws.Cells.Item(141, 11).Value = 280
$ws.Cells.Item(141, 12).Value = 280
$ws.Cells.Item(141, 13).Value = 280
$ws.Cells.Item(141, 16).Value = 280
$ws.Cells.Item(142, 4).Value = 44232
$ws.Cells.Item(142, 8).Value = "Camote"
$ws.Cells.Item(142, 9).Value = "1a nueva(o)"
$ws.Cells.Item(142, 10).Value = 800
$ws.Cells.Item(142, 11).Value = 300
$ws.Cells.Item(142, 12).Value = 300
$ws.Cells.Item(142, 13).Value = 300
$ws.Cells.Item(142, 16).Value = 300
$ws.Cells.Item(143, 4).Value = 44327
$ws.Cells.Item(144, 4).Value = 44510
$ws.Cells.Item(144, 8).Value = "Paine"
$ws.Cells.Item(144, 9).Value = "1a (guarda)"
$ws.Cells.Item(144, 10).Value = 2500
$ws.Cells.Item(144, 11).Value = 80
$ws.Cells.Item(144, 12).Value = 80
$ws.Cells.Item(144, 13).Value = 80
$ws.Cells.Item(144, 16).Value = 80
$ws.Cells.Item(145, 4).Value = 44468
$ws.Cells.Item(145, 8).Value = "Paine"
$ws.Cells.Item(145, 10).Value = 1500
$ws.Cells.Item(145, 11).Value = 150
$ws.Cells.Item(145, 12).Value = 150
$ws.Cells.Item(145, 13).Value = 150
$ws.Cells.Item(145, 16).Value = 150
$ws.Cells.Item(146, 4).Value = 44391
$ws.Cells.Item(146, 9).Value = "1a (guarda)"
$ws.Cells.Item(146, 10).Value = 900
$ws.Cells.Item(147, 4).Value = 44236
$ws.Cells.Item(147, 9).Value = "1a nueva(o)"
$ws.Cells.Item(147, 10).Value = 900
$ws.Cells.Item(147, 11).Value = 300
$ws.Cells.Item(147, 12).Value = 300
$ws.Cells.Item(147, 13).Value = 300
$ws.Cells.Item(147, 16).Value = 300
$ws.Cells.Item(148, 4).Value = 44389
$ws.Cells.Item(148, 10).Value = 900
$ws.Cells.Item(148, 11).Value = 300
$ws.Cells.Item(148, 12).Value = 300
$ws.Cells.Item(148, 13).Value = 300
$ws.Cells.Item(148, 16).Value = 300
$ws.Cells.Item(149, 4).Value = 44251
$ws.Cells.Item(149, 11).Value = 300
$ws.Cells.Item(149, 12).Value = 300
$ws.Cells.Item(149, 13).Value = 300
$ws.Cells.Item(149, 16).Value = 300
$ws.Cells.Item(150, 4).Value = 44330
$ws.Cells.Item(150, 9).Value = "1a (guarda)"
$ws.Cells.Item(150, 10).Value = 800
$ws.Cells.Item(150, 11).Value = 280
$ws.Cells.Item(150, 12).Value = 280
$ws.Cells.Item(150, 13).Value = 280
$ws.Cells.Item(150, 16).Value = 280
$ws.Cells.Item(151, 4).Value = 44432
$ws.Cells.Item(151, 9).Value = "1a (guarda)"
$ws.Cells.Item(151, 10).Value = 500
$ws.Cells.Item(151, 11).Value = 650
$ws.Cells.Item(151, 12).Value = 650
$ws.Cells.Item(151, 13).Value = 650
$ws.Cells.Item(151, 16).Value = 650
$ws.Cells.Item(152, 4).Value = 44181
$ws.Cells.Item(152, 8).Value = "Camote"
$ws.Cells.Item(152, 9).Value = "1a nueva(o)"
$ws.Cells.Item(152, 10).Value = 800
$ws.Cells.Item(152, 11).Value = 650
$ws.Cells.Item(152, 12).Value = 650
$ws.Cells.Item(152, 13).Value = 650
$ws.Cells.Item(152, 16).Value = 650
$ws.Cells.Item(153, 4).Value = 44194
$ws.Cells.Item(153, 9).Value = "1a nueva(o)"
$ws.Cells.Item(153, 10).Value = 1500
$ws.Cells.Item(153, 11).Value = 400
$ws.Cells.Item(153, 12).Value = 400
$ws.Cells.Item(153, 13).Value = 400
$ws.Cells.Item(153, 16).Value = 400
$ws.Cells.Item(154, 4).Value = 44271
$ws.Cells.Item(154, 9).Value = "1a (cosecha)"
$ws.Cells.Item(154, 10).Value = 1000
$ws.Cells.Item(154, 11).Value = 230
$ws.Cells.Item(154, 12).Value = 250
$ws.Cells.Item(154, 13).Value = 240
$ws.Cells.Item(154, 16).Value = 240
$ws.Cells.Item(155, 1).Value = 5
$ws.Cells.Item(155, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(155, 3).Value = "Maule"
$ws.Cells.Item(155, 4).Value = 44271
$ws.Cells.Item(155, 5).Value = 7
$ws.Cells.Item(155, 6).Value = 100112045
$ws.Cells.Item(155, 7).Value = "Zapallo"
$ws.Cells.Item(155, 8).Value = "Paine"
$ws.Cells.Item(155, 9).Value = "1a (cosecha)"
$ws.Cells.Item(155, 10).Value = 1200
$ws.Cells.Item(155, 11).Value = 150
$ws.Cells.Item(155, 12).Value = 150
$ws.Cells.Item(155, 13).Value = 150
$ws.Cells.Item(155, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(155, 15).Value = "Región del Maule"
$ws.Cells.Item(155, 16).Value = 150
$ws.Cells.Item(155, 17).Value = 1
$ws.Cells.Item(155, 18).Value = "Hortaliza"
$ws.Cells.Item(155, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(156, 1).Value = 5
$ws.Cells.Item(156, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(156, 3).Value = "Maule"
$ws.Cells.Item(156, 4).Value = 44400
$ws.Cells.Item(156, 5).Value = 7
$ws.Cells.Item(156, 6).Value = 100112045
$ws.Cells.Item(156, 7).Value = "Zapallo"
$ws.Cells.Item(156, 8).Value = "Camote"
$ws.Cells.Item(156, 9).Value = "1a (guarda)"
$ws.Cells.Item(156, 10).Value = 800
$ws.Cells.Item(156, 11).Value = 350
$ws.Cells.Item(156, 12).Value = 350
$ws.Cells.Item(156, 13).Value = 350
$ws.Cells.Item(156, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(156, 15).Value = "Región del Maule"
$ws.Cells.Item(156, 16).Value = 350
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"
$ws.Cells.Item(156, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(157, 1).Value = 5
$ws.Cells.Item(157, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(157, 3).Value = "Maule"
$ws.Cells.Item(157, 4).Value = 44201
$ws.Cells.Item(157, 5).Value = 7
$ws.Cells.Item(157, 6).Value = 100112045
$ws.Cells.Item(157, 7).Value = "Zapallo"
$ws.Cells.Item(157, 8).Value = "Camote"
$ws.Cells.Item(157, 9).Value = "1a nueva(o)"
$ws.Cells.Item(157, 10).Value = 800
$ws.Cells.Item(157, 11).Value = 350
$ws.Cells.Item(157, 12).Value = 350
$ws.Cells.Item(157, 13).Value = 350
$ws.Cells.Item(157, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(157, 15).Value = "Región del Maule"
$ws.Cells.Item(157, 16).Value = 350
$ws.Cells.Item(157, 17).Value = 1
$ws.Cells.Item(157, 18).Value = "Hortaliza"
$ws.Cells.Item(157, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
